# The document carries the Pearson / BTEC logo artwork as inline
# pictures inside its header and footers. Each inline picture's
# internal "name" metadata (InlineShape.Name) was swapped:
#   - the two Pearson logo pictures ("image2.png") became "image1.png"
#   - the BTEC logo picture ("image1.jpg") became "image2.jpg"
# Locate each picture via its description (AlternativeText), which
# identifies which logo it is, then rename it.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    # Footers: Pearson logo pictures -> rename to image1.png
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }

    # Headers: BTEC logo picture -> rename to image2.jpg
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }
}
